# Applies the changes from the commit "Funcao Gaussiana corrigida e testes"
# - Adds new test data to row 13 (columns O:W) on sheet "Mapa2a"
# - Updates selected cell on that sheet to V25

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapa2a")

# Fill in new test row (row 13, columns O:W) with the Gaussian detector test values
$ws.Range("O13").Value = 0.01
$ws.Range("P13").Value = 0.2
$ws.Range("Q13").Value = 0.04
$ws.Range("R13").Value = "0,1-0,7"
$ws.Range("S13").Value = "0-3"
$ws.Range("T13").Value = "Sim"
$ws.Range("U13").Value = "Não"
$ws.Range("V13").Value = 1
$ws.Range("W13").Value = 11

# Update the selected / active cell on the sheet
$ws.Activate()
$ws.Range("V25").Select()
